$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.724.90'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '3.490.87'
$ws.Range("E3").Value = '  +4.96%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = "'249.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("D6").Value = "'659.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.18%  '

$ws.Range("D7").Value = "'1.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.98%  '

$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("E10").Value = '  +1.93%  '

$ws.Range("D11").Value = '3.489.29'
$ws.Range("E11").Value = '  +5.01%  '

$ws.Range("D12").Value = "'44.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +10.57%  '

$ws.Range("E13").Value = '  +0.85%  '

$ws.Range("D14").Value = '97.488.62'
$ws.Range("E14").Value = '  +0.42%  '

$ws.Range("D15").Value = "'6.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.67%  '

$ws.Range("D16").Value = '4.147.77'
$ws.Range("E16").Value = '  +5.03%  '

$ws.Range("E17").Value = '  +2.02%  '

$ws.Range("D18").Value = "'8.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("D19").Value = '3.484.91'
$ws.Range("E19").Value = '  +4.68%  '

$ws.Range("D20").Value = "'18.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.02%  '

$ws.Range("D21").Value = "'12.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +14.49%  '

$ws.Range("D22").Value = "'0.500"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.00%  '

$ws.Range("D23").Value = "'519.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.93%  '

$ws.Range("E24").Value = '  +1.77%  '

$ws.Range("E25").Value = '  +0.51%  '

$ws.Range("D26").Value = "'6.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.17%  '

$ws.Range("D27").Value = "'96.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.31%  '

$ws.Range("D28").Value = "'12.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.71%  '

$ws.Range("D29").Value = '3.672.37'
$ws.Range("E29").Value = '  +4.82%  '

$ws.Range("D30").Value = "'12.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +14.14%  '

$ws.Range("D31").Value = "'2.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +15.94%  '

$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.62%  '

$ws.Range("E33").Value = '  -2.10%  '

$ws.Range("D34").Value = "'0.187"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("D35").Value = "'0.595"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.39%  '

$ws.Range("D36").Value = "'31.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.16%  '

$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("D38").Value = "'7.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.56%  '

$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("D40").Value = "'0.155"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.92%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = "'519.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.17%  '

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").Value = "'0.915"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.51%  '

$ws.Range("D44").Value = "'24.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.85%  '

$ws.Range("E45").Value = '  +4.71%  '

$ws.Range("E46").Value = '  +3.70%  '

$ws.Range("E47").Value = '  +4.07%  '

$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = "'3.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.37%  '

$ws.Range("B49").Value = 'MantraDAO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D49").Value = "'3.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.95%  '

$ws.Range("E50").Value = '  +12.51%  '

$ws.Range("E51").Value = '  -0.54%  '
